$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.602.14'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.268.51'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.44%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '120.79'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.42%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.67%  '

$ws.Range('E7').Value = '  +4.12%  '

$ws.Range('E8').Value = '  +0.26%  '

$ws.Range('E9').Value = '  +2.29%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.72'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.10%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0944'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.45'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +6.74%  '

$ws.Range('E13').Value = '  -1.83%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.87'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.26%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.914'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.82%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.609.33'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.40%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.264.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.598.62'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.91%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000109'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.45%  '

$ws.Range('E20').Value = '  +2.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.89%  '

$ws.Range('E22').Value = '  -5.07%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.06'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.23%  '

$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.73%  '

$ws.Range('E27').Value = '  +1.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.12'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.48%  '

$ws.Range('E29').Value = '  +0.11%  '

$ws.Range('E30').Value = '  +0.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.99'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.72%  '

$ws.Range('E32').Value = '  +0.96%  '

$ws.Range('E33').Value = '  +1.09%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.71'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.14%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.50'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +16.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.130'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.94%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0381'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.66%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.72'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.109'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.99%  '

$ws.Range('E40').Value = '  -3.83%  '

$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.77'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.18%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.242'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.34'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -5.32%  '

$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('E45').Value = '  -0.71%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.70'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -7.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '76.39'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +38.23%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.671'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +19.87%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.58'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.78%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.88'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.45%  '
